$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-05-09 Friday" "2025-05-10 Saturday"

Replace-Text "43×90=" "22×93="
Replace-Text "43×14=" "69×76="
Replace-Text "80×56=" "65×64="
Replace-Text "68×65=" "58×29="
Replace-Text "99×21=" "74×55="
Replace-Text "27×51=" "66×89="
Replace-Text "41×64=" "37×40="
Replace-Text "72×57=" "90×72="
Replace-Text "74×16=" "27×58="
Replace-Text "53×39=" "78×78="
Replace-Text "65×93=" "61×97="
Replace-Text "22×70=" "88×16="
Replace-Text "33×46=" "76×11="
Replace-Text "18×46=" "72×35="
Replace-Text "43×82=" "88×40="
Replace-Text "14×17=" "32×30="
Replace-Text "55×20=" "94×84="
Replace-Text "79×79=" "22×51="
Replace-Text "27×11=" "62×37="
Replace-Text "61×74=" "24×81="
Replace-Text "44×53=" "72×18="
Replace-Text "48×68=" "12×17="
Replace-Text "59×73=" "59×79="
Replace-Text "59×75=" "75×47="
Replace-Text "89×68=" "50×16="
